# Apply the edit described by the diff:
# - "2016-09-01" sheet loses tabSelected, selection moves to A6:G6
# - "2017-01-01" sheet: selection moves to A6:G6; row 6 gets D6=1 and
#   E6/F6/G6 text updated (department/city/job_title)
# - "2017-03-12" sheet becomes the active/tabSelected sheet, selection moves
#   to A6:G6; row 6 gets the same D6/E6/F6/G6 update
# - workbook active tab becomes "2017-03-12" (index 5)

$wb = $excel.ActiveWorkbook

# --- Update "2017-01-01" sheet (row 6 data + selection) ---
$wsJan = $wb.Worksheets.Item("2017-01-01")
$wsJan.Range("D6").Value = 1
$wsJan.Range("E6").Value = "MARK/ENG/PRI/AME"
$wsJan.Range("F6").Value = "New York"
$wsJan.Range("G6").Value = "VIE"

# --- Update "2016-09-01" sheet (no data change, only selection/tab state) ---
$wsSep = $wb.Worksheets.Item("2016-09-01")
$wsSep.Activate()
$wsSep.Range("A6:G6").Select()

# --- Update "2017-03-12" sheet (row 6 data + becomes active/selected sheet) ---
$wsMar = $wb.Worksheets.Item("2017-03-12")
$wsMar.Range("D6").Value = 1
$wsMar.Range("E6").Value = "MARK/ENG/PRI/AME"
$wsMar.Range("F6").Value = "New York"
$wsMar.Range("G6").Value = "VIE"

# Re-apply the A6:G6 selection on "2017-01-01" as its own last-used view
$wsJan.Activate()
$wsJan.Range("A6:G6").Select()

# Finally activate "2017-03-12" so it becomes the tabSelected / active tab
$wsMar.Activate()
$wsMar.Range("A6:G6").Select()
